# Add feeding history to group (Containers) and individual (Details) reports

$wb = $excel.ActiveWorkbook

$wsContainers = $wb.Worksheets.Item("Containers")
$wsDetails = $wb.Worksheets.Item("Details")

# --- Details sheet: update the stale selection (no longer the active tab) ---
$wsDetails.Activate()
$wsDetails.Range("O1:T5").Select()

# --- Containers sheet: add the new "Feedings" block in columns O:T ---
$wsContainers.Activate()

# header label above the new table
$wsContainers.Range("O1").Value = "Feedings"

# new table header row (row 4), matching the styling of the existing tables
$wsContainers.Range("O4").Value = "Feed"
$wsContainers.Range("O4").Style = $wsContainers.Range("C4").Style

$wsContainers.Range("P4").Value = "Date Started"
$wsContainers.Range("P4").Style = $wsContainers.Range("C4").Style

$wsContainers.Range("Q4").Value = "Container"
$wsContainers.Range("Q4").Style = $wsContainers.Range("C4").Style

$wsContainers.Range("R4").Value = "Frequency"
$wsContainers.Range("R4").Style = $wsContainers.Range("J4").Style

$wsContainers.Range("S4").Value = "Method"
$wsContainers.Range("S4").Style = $wsContainers.Range("J4").Style

$wsContainers.Range("T4").Value = "Comments"
$wsContainers.Range("T4").Style = $wsContainers.Range("J4").Style

# widen the new columns to fit the new headers
$wsContainers.Columns.Item(15).ColumnWidth = 14.736979166666666
$wsContainers.Columns.Item(16).ColumnWidth = 14.022135416666666
$wsContainers.Columns.Item(17).ColumnWidth = 13.166666666666666
$wsContainers.Columns.Item(18).ColumnWidth = 18.307291666666668
$wsContainers.Columns.Item(19).ColumnWidth = 15.166666666666666
$wsContainers.Columns.Item(20).ColumnWidth = 26.877604166666668

# Containers becomes the active tab/sheet, with the new selection
$wsContainers.Range("R10").Select()
